# Apply "Natmi following Dr Hou advice" update: replace the 3 existing
# ECs/FAPs/sCs -> sCs data rows with a full 3x3 Cartesian grid (9 rows)
# of sending-cluster x target-cluster combinations for Ccl11-Ccr2,
# using the recomputed statistics for each pairing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @(2, @("ECs", "Ccl11", "Ccr2", "ECs", 3, 1, 0.7503183333333333, 2.250955, 0.001871730106429624, 0.001871730106429624, 2, 0.6666666666666666, 9.235433666666667, 27.706301, 0.9849159412561933, 0.9849159412561931, 6.929515196383889, 62.365636767455, 0.001843496819551688, 0.001843496819551688)),
  @(3, @("ECs", "Ccl11", "Ccr2", "FAPs", 3, 1, 0.7503183333333333, 2.250955, 0.001871730106429624, 0.001871730106429624, 2, 0.6666666666666666, 0.08345333333333332, 0.25036, 0.008899908907107467, 0.008899908907107465, 0.06261656597777777, 0.5635490937999998, [double]"1.665822744591422e-05", [double]"1.665822744591422e-05")),
  @(4, @("ECs", "Ccl11", "Ccr2", "sCs", 3, 1, 0.7503183333333333, 2.250955, 0.001871730106429624, 0.001871730106429624, 3, 1, 0.057988, 0.173964, 0.006184149836699327, 0.006184149836699326, 0.04350945951333333, 0.39158513562, [double]"1.157505943202198e-05", [double]"1.157505943202197e-05")),
  @(5, @("FAPs", "Ccl11", "Ccr2", "ECs", 3, 1, 383.1307676666667, 1149.392303, 0.9557508602453543, 0.9557508602453542, 2, 0.6666666666666666, 9.235433666666667, 27.706301, 0.9849159412561933, 0.9849159412561931, 3538.378790444579, 31845.40911400121, 0.9413342581249695, 0.9413342581249693)),
  @(6, @("FAPs", "Ccl11", "Ccr2", "FAPs", 3, 1, 383.1307676666667, 1149.392303, 0.9557508602453543, 0.9557508602453542, 2, 0.6666666666666666, 0.08345333333333332, 0.25036, 0.008899908907107467, 0.008899908907107465, 31.97353966434222, 287.76185697908, 0.008506095594073252, 0.008506095594073249)),
  @(7, @("FAPs", "Ccl11", "Ccr2", "sCs", 3, 1, 383.1307676666667, 1149.392303, 0.9557508602453543, 0.9557508602453542, 3, 1, 0.057988, 0.173964, 0.006184149836699327, 0.006184149836699326, 22.21698695545467, 199.952882599092, 0.005910506526311548, 0.005910506526311548)),
  @(8, @("sCs", "Ccl11", "Ccr2", "ECs", 3, 1, 16.98778433333333, 50.963353, 0.04237740964821621, 0.0423774096482162, 2, 0.6666666666666666, 9.235433666666667, 27.706301, 0.9849159412561933, 0.9849159412561931, 156.8895553541392, 1412.005998187253, 0.04173818631167216, 0.04173818631167214)),
  @(9, @("sCs", "Ccl11", "Ccr2", "FAPs", 3, 1, 16.98778433333333, 50.963353, 0.04237740964821621, 0.0423774096482162, 2, 0.6666666666666666, 0.08345333333333332, 0.25036, 0.008899908907107467, 0.008899908907107465, 1.417687228564444, 12.75918505708, 0.0003771550855883013, 0.0003771550855883012)),
  @(10, @("sCs", "Ccl11", "Ccr2", "sCs", 3, 1, 16.98778433333333, 50.963353, 0.04237740964821621, 0.0423774096482162, 3, 1, 0.057988, 0.173964, 0.006184149836699327, 0.006184149836699326, 0.9850876379213335, 8.865788741292, 0.0002620682509557567, 0.0002620682509557567))
)

foreach ($pair in $rows) {
  $r = $pair[0]
  $vals = $pair[1]
  for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
  }
}
